$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: insert a new "Menu List for Navigation Bar" bullet (numId 3)
# right after the "completedItem" bullet.
# ---------------------------------------------------------------------
$pCompleted = $d.Paragraphs.Item(22)
$rCompleted = $pCompleted.Range
$rCompleted.Start = $rCompleted.End - 1
$rCompleted.Collapse(0)
$rCompleted.InsertParagraphAfter()
$pMenu = $d.Paragraphs.Item(23)
$pMenu.Range.Text = "Menu List for Navigation Bar"

# ---------------------------------------------------------------------
# Change 2: split the "Navigation Bar " paragraph into three runs:
#   "Navigation Bar" + "- " + " "
# ---------------------------------------------------------------------
$pNav = $d.Paragraphs.Item(26)
$rNav = $pNav.Range
$rNav.End = $rNav.End - 1
$rNav.Text = "Navigation Bar"

$rNavEnd = $pNav.Range
$rNavEnd.Start = $rNavEnd.End - 1
$rNavEnd.Collapse(0)
$rNavEnd.InsertParagraphAfter()
$pDash = $d.Paragraphs.Item(27)
$pDash.Range.Text = "- "

$rDashEnd = $pDash.Range
$rDashEnd.Start = $rDashEnd.End - 1
$rDashEnd.Collapse(0)
$rDashEnd.InsertParagraphAfter()
$pSpace = $d.Paragraphs.Item(28)
$pSpace.Range.Text = " "

# merge the three paragraphs back together (removes the two inserted
# paragraph marks) so the three sentences become three runs of one
# paragraph
$mark1 = $d.Paragraphs.Item(26).Range
$m1 = $mark1.End - 1
$d.Range($m1, $m1 + 1).Delete()

$mark2 = $d.Paragraphs.Item(26).Range
$m2 = $mark2.End - 1
$d.Range($m2, $m2 + 1).Delete()

# ---------------------------------------------------------------------
# Change 3: "Comments field" -> "Comments"[gramStart/gramEnd] + " field"
# + "/ Message Field"
# Change 4: after that paragraph, add a blank paragraph then three new
# lines describing the app structure.
# Both changes are applied together via a single InsertXML call so the
# inserted blank paragraph stays a truly empty <w:p/> (no stray <w:r/>).
# ---------------------------------------------------------------------
$pComments = $d.Paragraphs.Item(32)
$rComments = $pComments.Range
$rComments.End = $rComments.End - 1

$endash = [char]0x2013
$commentsXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>Comments</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> field</w:t></w:r>
<w:r><w:t>/ Message Field</w:t></w:r>
</w:p>
<w:p/>
<w:p><w:r><w:t xml:space="preserve">App -&gt; Nav Bar -&gt; Menu List, CSS </w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">             To Do List $endash To Do List, To Do Form</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">             Contact Form $endash Contact Form, CSS</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$null = $rComments.InsertXML($commentsXml)
